# Update the "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets. Both sheets carry the same exhibition rows, so the same
# row -> new value mapping is applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 110
    3  = 402
    4  = 11773
    5  = 915
    6  = 121
    7  = 21
    9  = 150
    11 = 26
    17 = 1395
    19 = 910
    20 = 112
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
